# "Generate Report for Handoff" - refresh the localization-status report:
#   - Status moves from "Handed back: in sync with en-US" to "In Translation"
#     (and the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are refreshed) for all three sheets.
#   - zh-cn / de-de sheets get a new "Error Detail" message because the
#     handback file used is stale compared to the latest handoff.
#   - The now-unused wide "Status" / "Error Detail" columns are narrowed /
#     widened to their new natural size.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1b5de1721502e560b6602d8f4e47f0de9cc713a/e2e/24ccb9c7-e03d-4498-af8f-4682dacd4df7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/243ee3f6d4d480fc7bb30bd8c2786872bdce81b5/e2e/24ccb9c7-e03d-4498-af8f-4682dacd4df7.md."

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2017-02-09 13:47:09"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2017-02-09 13:46:52"
$wsZhCn.Range("R2").Value = $errorDetail
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsZhCn.Columns.Item(18).ColumnWidth = 39.166666666666664

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("R2").Value = $errorDetail
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(18).ColumnWidth = 39.166666666666664
